$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59, col A: tiny precision update to the stored timestamp
$ws.Cells.Item(59, 1).Value = 44372.76793294908

# New row 60: newly retrieved data point
$ws.Cells.Item(60, 1).Value = 44373.76802768324
$ws.Cells.Item(60, 1).NumberFormat = $ws.Cells.Item(59, 1).NumberFormat

$ws.Cells.Item(60, 2).Value = 78959
$ws.Cells.Item(60, 3).Value = 66353
$ws.Cells.Item(60, 4).Value = 3631
$ws.Cells.Item(60, 5).Value = 2149
$ws.Cells.Item(60, 6).Value = 1529
$ws.Cells.Item(60, 7).Value = 20917
$ws.Cells.Item(60, 8).Value = 1546
$ws.Cells.Item(60, 9).Value = 887
$ws.Cells.Item(60, 10).Value = 193
